$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7864
$ws.Range("C2").Value = 131.0666666666667
$ws.Range("D2").Value = 0.1330889652511585
$ws.Range("E2").Value = 8.981193069482002
$ws.Range("F2").Value = 26.5306965659216
